$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update URL value (row 2)
$ws.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/PreferredAllowedReason"

# Update Date value (row 8)
$ws.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# Insert a new row before row 11 ("Description" row) for the new "Jurisdiction" property
$ws.Rows("11:11").Insert()

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "'"
